# Hotfix: Tue Nov 19 17:53:55 RTZ 2024
#
# The "SQL" sheet stored each row as (Id, SqlCommand, RussianDescription).
# The new layout is (Id, RussianDescription, SqlCommand, Timestamp, <spare>):
#   - column B (old SQL command) and column C (old Russian description)
#   swap places, so B now holds the human-readable description and C the
#   SQL text;
#   - a new column D is added with a creation timestamp for every row;
#   - a new (currently empty) column E is appended as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SQL")

$timestamp = "2024-11-19 13:58:50"
$lastRow = 13

for ($r = 1; $r -le $lastRow; $r++) {
    $oldCommand = $ws.Cells.Item($r, 2).Value()
    $oldDescription = $ws.Cells.Item($r, 3).Value()

    # Swap: description moves to B, SQL command moves to C.
    $ws.Cells.Item($r, 2).Value = $oldDescription
    $ws.Cells.Item($r, 3).Value = $oldCommand

    # New timestamp column.
    $ws.Cells.Item($r, 4).Value = $timestamp

    # New trailing column E, left blank but materialised as a real cell
    # (touching a formatting property forces Excel to keep the cell
    # without assigning it any content or a non-default style).
    $ws.Cells.Item($r, 5).Font.Bold = $ws.Cells.Item($r, 5).Font.Bold
}
